# Apply the 2024-07-19 cryptos-list refresh (prices / 1h-volume deltas,
# plus the Monero <-> ImmutableX row-order swap at rows 36-37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.029.26'
$ws.Range("E2").Value = '  -0.89%  '
# Row 3
$ws.Range("D3").Value = '3.402.08'
$ws.Range("E3").Value = '  -1.15%  '
# Row 4
$ws.Range("E4").Value = '  -0.04%  '
# Row 5
$ws.Range("D5").Value = '''570.60'
$ws.Range("E5").Value = '  -0.31%  '
# Row 6
$ws.Range("D6").Value = '''162.31'
$ws.Range("E6").Value = '  +2.38%  '
# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.01%  '
# Row 8
$ws.Range("D8").Value = '3.402.81'
$ws.Range("E8").Value = '  -1.12%  '
# Row 9
$ws.Range("D9").Value = '''0.547'
$ws.Range("E9").Value = '  -4.26%  '
# Row 10
$ws.Range("D10").Value = '''7.31'
$ws.Range("E10").Value = '  +1.48%  '
# Row 11
$ws.Range("D11").Value = '''0.119'
$ws.Range("E11").Value = '  -1.65%  '
# Row 12
$ws.Range("E12").Value = '  -4.17%  '
# Row 13
$ws.Range("D13").Value = '3.990.46'
$ws.Range("E13").Value = '  -1.18%  '
# Row 14
$ws.Range("E14").Value = '  +0.74%  '
# Row 15
$ws.Range("D15").Value = '''26.82'
$ws.Range("E15").Value = '  -2.27%  '
# Row 16
$ws.Range("D16").Value = '''0.0000172'
$ws.Range("E16").Value = '  -0.51%  '
# Row 17
$ws.Range("D17").Value = '64.075.01'
$ws.Range("E17").Value = '  -0.93%  '
# Row 18
$ws.Range("D18").Value = '3.388.36'
$ws.Range("E18").Value = '  -1.49%  '
# Row 19
$ws.Range("D19").Value = '''6.10'
$ws.Range("E19").Value = '  -0.57%  '
# Row 20
$ws.Range("D20").Value = '''13.47'
$ws.Range("E20").Value = '  -1.38%  '
# Row 21
$ws.Range("D21").Value = '''371.59'
$ws.Range("E21").Value = '  -1.31%  '
# Row 22
$ws.Range("D22").Value = '''7.77'
$ws.Range("E22").Value = '  -1.20%  '
# Row 23
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.31%  '
# Row 24
$ws.Range("D24").Value = '''70.04'
$ws.Range("E24").Value = '  -2.75%  '
# Row 25
$ws.Range("D25").Value = '''0.510'
$ws.Range("E25").Value = '  -3.95%  '
# Row 26
$ws.Range("E26").Value = '  -3.96%  '
# Row 27
$ws.Range("D27").Value = '''9.46'
$ws.Range("E27").Value = '  -4.24%  '
# Row 28
$ws.Range("E28").Value = '  -0.43%  '
# Row 29
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.21%  '
# Row 30
$ws.Range("D30").Value = '''6.05'
$ws.Range("E30").Value = '  +0.38%  '
# Row 31
$ws.Range("D31").Value = '''1.39'
$ws.Range("E31").Value = '  -2.90%  '
# Row 32
$ws.Range("E32").Value = '  -0.71%  '
# Row 33
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  +0.06%  '
# Row 34
$ws.Range("D34").Value = '''22.69'
$ws.Range("E34").Value = '  -1.84%  '
# Row 35
$ws.Range("D35").Value = '''6.96'
$ws.Range("E35").Value = '  +0.12%  '
# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''1.47'
$ws.Range("E36").Value = '  -5.58%  '
# Row 37
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''159.51'
$ws.Range("E37").Value = '  -0.84%  '
# Row 38
$ws.Range("D38").Value = '''0.859'
$ws.Range("E38").Value = '  +8.99%  '
# Row 39
$ws.Range("D39").Value = '''1.80'
$ws.Range("E39").Value = '  -3.94%  '
# Row 40
$ws.Range("D40").Value = '''25.79'
$ws.Range("E40").Value = '  -0.93%  '
# Row 41
$ws.Range("D41").Value = '''0.0719'
$ws.Range("E41").Value = '  -3.27%  '
# Row 42
$ws.Range("D42").Value = '''42.58'
$ws.Range("E42").Value = '  -0.98%  '
# Row 43
$ws.Range("D43").Value = '2.730.43'
$ws.Range("E43").Value = '  -5.27%  '
# Row 44
$ws.Range("D44").Value = '''6.41'
$ws.Range("E44").Value = '  -0.34%  '
# Row 45
$ws.Range("D45").Value = '''25.78'
$ws.Range("E45").Value = '  -0.52%  '
# Row 46
$ws.Range("D46").Value = '''4.33'
$ws.Range("E46").Value = '  -3.62%  '
# Row 47
$ws.Range("E47").Value = '  -1.98%  '
# Row 48
$ws.Range("D48").Value = '''2.39'
$ws.Range("E48").Value = '  +1.11%  '
# Row 49
$ws.Range("D49").Value = '''328.44'
$ws.Range("E49").Value = '  +2.34%  '
# Row 50
$ws.Range("D50").Value = '''1.04'
$ws.Range("E50").Value = '  -3.54%  '
# Row 51
$ws.Range("E51").Value = '  -1.81%  '
